$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44319, 2, 14, 199.1465149359886),
    @(44320, 0, 13, 184.9217638691323),
    @(44321, 1, 14, 199.1465149359886)
)

$row = 245
foreach ($r in $data) {
    $ws.Cells.Item($row - 1, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]

    $row++
}
